{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block (and\n// the blank paragraph immediately preceding it) that followed the last\n// bibliography line (\"Rio de Janeiro: Elsevier Editora, 2007.\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph (last bibliography entry) by its exact text.\nconst anchorText = \"Rio de Janeiro: Elsevier Editora, 2007.\";\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find anchor paragraph: \" + anchorText);\n}\n\n// The three paragraphs that must go: the blank paragraph right after the\n// anchor, the \"Ver no Jupiter ...\" line, and the \"\u00a9 2020 ...\" line.\nconst expectedRemoved = [\n  \"\",\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst toDelete = [];\nfor (let k = 0; k < expectedRemoved.length; k++) {\n  const idx = anchorIndex + 1 + k;\n  if (idx >= items.length || items[idx].text !== expectedRemoved[k]) {\n    throw new Error(\"Unexpected document shape near anchor; refusing to delete.\");\n  }\n  toDelete.push(items[idx]);\n}\n\n// Delete from the end backwards so earlier proxies in the batch stay valid.\nfor (let k = toDelete.length - 1; k >= 0; k--) {\n  toDelete[k].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"(c) 2020 ...\" footer block (and\n# the blank paragraph immediately preceding it) that followed the last\n# bibliography line (\"Rio de Janeiro: Elsevier Editora, 2007.\").\n\n$d = $word.ActiveDocument\n\nfunction Get-ParaText($para) {\n    # Paragraph.Range.Text includes the trailing paragraph/cell mark; strip it\n    # so comparisons match the visible text exactly.\n    return $para.Range.Text.TrimEnd([char]13, [char]7)\n}\n\n$count = $d.Paragraphs.Count\n\n$anchorText = \"Rio de Janeiro: Elsevier Editora, 2007.\"\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    if ((Get-ParaText $d.Paragraphs($i)) -eq $anchorText) {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not find anchor paragraph: $anchorText\"\n}\n\n# The three paragraphs that must go: the blank paragraph right after the\n# anchor, the \"Ver no Jupiter ...\" line, and the \"(c) 2020 ...\" line.\n$expectedRemoved = @(\n    \"\",\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    [char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n)\n\n# Validate the paragraphs are exactly what we expect before touching anything.\nfor ($k = 0; $k -lt $expectedRemoved.Length; $k++) {\n    $idx = $anchorIndex + 1 + $k\n    if ($idx -gt $d.Paragraphs.Count) {\n        throw \"Unexpected document shape near anchor; refusing to delete.\"\n    }\n    $actual = Get-ParaText $d.Paragraphs($idx)\n    if ($actual -ne $expectedRemoved[$k]) {\n        throw \"Unexpected document shape near anchor; refusing to delete.\"\n    }\n}\n\n# Delete from the end backwards so earlier (lower) indices stay valid.\nfor ($k = $expectedRemoved.Length - 1; $k -ge 0; $k--) {\n    $idx = $anchorIndex + 1 + $k\n    $d.Paragraphs($idx).Range.Delete()\n}\n"}
